# Update Cd80-Cd28 LR-pairs sheet with new TPM-derived values.
# The new data collapses the old 8 combinations (rows 2-9) down to 4 rows
# (rows 2-5), one per Sending cluster (ECs, FAPs, MuSCs, Resolving-Mac),
# all now paired against "Resolving-Mac" as the Target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (old extra rows no longer present)
$ws.Range("A6:T9").EntireRow.Delete() | Out-Null

# Row 2: ECs -> Resolving-Mac
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd80"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7715313333333333
$ws.Range("H2").Value = 2.314594
$ws.Range("I2").Value = 0.05172308417778351
$ws.Range("J2").Value = 0.05172308417778351
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.961972333333333
$ws.Range("N2").Value = 14.885917
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3.828317130299777
$ws.Range("R2").Value = 34.454854172698
$ws.Range("S2").Value = 0.05172308417778351
$ws.Range("T2").Value = 0.05172308417778351

# Row 3: FAPs -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cd80"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.713472666666666
$ws.Range("H3").Value = 11.140418
$ws.Range("I3").Value = 0.2489493958723191
$ws.Range("J3").Value = 0.2489493958723191
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.961972333333333
$ws.Range("N3").Value = 14.885917
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 18.42614863258955
$ws.Range("R3").Value = 165.835337693306
$ws.Range("S3").Value = 0.2489493958723191
$ws.Range("T3").Value = 0.2489493958723191

# Row 4: MuSCs -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cd80"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.701035666666667
$ws.Range("H4").Value = 5.103107
$ws.Range("I4").Value = 0.1140366012048922
$ws.Range("J4").Value = 0.1140366012048922
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.961972333333333
$ws.Range("N4").Value = 14.885917
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 8.440491916013221
$ws.Range("R4").Value = 75.96442724411899
$ws.Range("S4").Value = 0.1140366012048922
$ws.Range("T4").Value = 0.1140366012048922

# Row 5: Resolving-Mac -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Cd80"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.730536666666668
$ws.Range("H5").Value = 26.19161
$ws.Range("I5").Value = 0.5852909187450052
$ws.Range("J5").Value = 0.5852909187450052
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.961972333333333
$ws.Range("N5").Value = 14.885917
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 43.32068139515222
$ws.Range("R5").Value = 389.88613255637
$ws.Range("S5").Value = 0.5852909187450052
$ws.Range("T5").Value = 0.5852909187450052
